# Regenerate save_data column G ("K") values for rows 2..38.
# These are the newly computed strikeout (K) counts replacing the old
# Strike# values, per the commit "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(9, 14, 6, 4, 9, 6, 11, 6, 3, 3, 6, 7, 4, 8, 7, 10, 0, 9, 6, 8, 4, 6, 7, 7, 6, 1, 9, 4, 8, 5, 11, 12, 4, 6, 1, 4, 3)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
